$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 61.04160633333334
$ws.Range("N2").Value = 183.124819
$ws.Range("O2").Value = 0.2043613460574534
$ws.Range("P2").Value = 0.2043613460574534
$ws.Range("Q2").Value = 34.03500913768089
$ws.Range("R2").Value = 306.3150822391281
$ws.Range("S2").Value = 0.2043613460574534
$ws.Range("T2").Value = 0.2043613460574534

# Row 3 updates
$ws.Range("O3").Value = 0.3559304658284363
$ws.Range("P3").Value = 0.3559304658284363
$ws.Range("S3").Value = 0.3559304658284363
$ws.Range("T3").Value = 0.3559304658284363

# Row 4 updates
$ws.Range("M4").Value = 131.3384093333333
$ws.Range("N4").Value = 394.015228
$ws.Range("O4").Value = 0.4397081881141102
$ws.Range("P4").Value = 0.4397081881141103
$ws.Range("Q4").Value = 73.23044445092621
$ws.Range("R4").Value = 659.074000058336
$ws.Range("S4").Value = 0.4397081881141102
$ws.Range("T4").Value = 0.4397081881141103
